$d = $word.ActiveDocument

# Remove the duplicated "αστερισμό του" before "Αστερισμός του Ηρακλή".
# This turns "... για τον αστερισμό του Αστερισμός του Ηρακλή: ..."
# into      "... για τον  Αστερισμός του Ηρακλή: ..." (note resulting double space).
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute(
    "αστερισμό του Αστερισμός του Ηρακλή",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    " Αστερισμός του Ηρακλή",
    2
)
